$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the tiny floating-point precision change on the existing last row (A16)
$ws.Range("A16").Value = 44329.77915936575

# Append the new data row (row 17)
$ws.Range("A17").Value = 44330.77832976371
$ws.Range("B17").Value = 74531
$ws.Range("C17").Value = 62615
$ws.Range("D17").Value = 3384
$ws.Range("E17").Value = 2116
$ws.Range("F17").Value = 1498
$ws.Range("G17").Value = 19409
$ws.Range("H17").Value = 1404
$ws.Range("I17").Value = 873
$ws.Range("J17").Value = 225

# Match the date-style formatting used for column A on the other rows
$ws.Range("A17").NumberFormat = $ws.Range("A16").NumberFormat
